$wb = $excel.ActiveWorkbook

# Remember the sheet that is active before we start, so we can restore the
# original active-tab / selection state once we're done adding sheets.
$originalActive = $wb.ActiveSheet

$src = $wb.Worksheets.Item("SpecsDataCalib2")

# --- Add "SpecsDataCalib3" as a copy of "SpecsDataCalib2" ---
$src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$calib3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$calib3.Name = "SpecsDataCalib3"
$calib3.Range("H2").Value = 0.1609394603523363
$calib3.Range("X2").Value = 0.5104112205648693
$calib3.Range("Y2").Value = 0.0331974777603167

# --- Add "SpecsDataCalib4" as a copy of "SpecsDataCalib2" ---
$src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$calib4 = $wb.Worksheets.Item($wb.Worksheets.Count)
$calib4.Name = "SpecsDataCalib4"
$calib4.Range("H2").Value = 0.1609394603523363
$calib4.Range("X2").Value = 0.5104112205648693
$calib4.Range("Y2").Value = 0.0331974777603167

# Restore the workbook's original active sheet/tab selection.
$originalActive.Activate()
